$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply every existing value in B2:B28 by 1000 (unit conversion, e.g. from
# thousands to full units).
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value()
    $cell.Value = $old * 1000
}

# Update view state: zoom + new selected cell on the worksheet.
$excel.ActiveWindow.Zoom = 85
$ws.Range("B2").Select()
